$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds numbers that are authored as text (shared strings), e.g.
# "79", "63", etc. Assigning a bare numeric-looking string would make Excel
# coerce the cell to a real number, so each value is entered with the
# classic leading-apostrophe text-entry prefix. That leaves a "quote
# prefix" style flag on the cell, so ClearFormats() is called right after
# to drop back to the default (unstyled) cell format -- matching how the
# sheet was originally authored -- while the stored value stays text.
$ws.Range("C1").Value = "'55"
$ws.Range("C1").ClearFormats()

$ws.Range("C2").Value = "'50"
$ws.Range("C2").ClearFormats()

$ws.Range("C3").Value = "'70"
$ws.Range("C3").ClearFormats()

$ws.Range("C4").Value = "'62"
$ws.Range("C4").ClearFormats()

$ws.Range("C5").Value = "'75"
$ws.Range("C5").ClearFormats()

$ws.Range("C6").Value = "'75"
$ws.Range("C6").ClearFormats()

$ws.Range("C7").Value = "'59"
$ws.Range("C7").ClearFormats()

$ws.Range("C8").Value = "'53"
$ws.Range("C8").ClearFormats()
